$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("power")

# ---------------------------------------------------------------------------
# 1) Steering controller row (row 5): add supply-voltage note + numeric 24V
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = "controller voltage supply  12-50V. motor input 0-50V. Wattage is not include motor wattage"
$ws.Range("D5").Value = 24
$ws.Rows(5).RowHeight = 45

# ---------------------------------------------------------------------------
# 2) Brake controller row (row 7): same controller + same supply note as
#    the steering controller above (reuses the shared string from step 1).
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "Faulhaber MC-5010S"
$ws.Range("C7").Value = "controller voltage supply  12-50V. motor input 0-50V. Wattage is not include motor wattage"
$ws.Rows(7).RowHeight = 45

# ---------------------------------------------------------------------------
# 3) Insert six new rows (16:21) for anti-roll servos, wing actuation and
#    cooling system entries, just above the "sensors" section.
# ---------------------------------------------------------------------------
$ws.Rows("16:21").Insert()

$ws.Range("A17").Value = "anti rool servos (4)"

$ws.Range("A16").Value = "anti rool servos controller (4) "
$ws.Rows(16).RowHeight = 30

# ---------------------------------------------------------------------------
# 4) Brake motor row (row 8): note + max wattage estimate
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = "estimate."
$ws.Range("F8").Value = 100

# ---------------------------------------------------------------------------
# 5) Steering motor row (row 6): add max wattage estimate
# ---------------------------------------------------------------------------
$ws.Range("F6").Value = 200

$ws.Range("A18").Value = "wing controller"
$ws.Range("B18").Value = "Faulhaber MC-5010S"
$ws.Range("B18").HorizontalAlignment = -4108
$ws.Range("B18").VerticalAlignment = -4108
$ws.Range("B18").WrapText = $true

$ws.Range("A19").Value = "wing motor"

$ws.Range("A20").Value = "colling pump"
$ws.Range("D20").Value = 24
$ws.Range("D20").HorizontalAlignment = -4108
$ws.Range("D20").WrapText = $true

$ws.Range("A21").Value = "cooling fan"
$ws.Range("C21").Value = "SPAL VA32-A101-62A 12V 3.4A-MAX"
$ws.Range("C21").HorizontalAlignment = -4108
$ws.Range("C21").WrapText = $true

# ---------------------------------------------------------------------------
# 6) Widen column A to fit the new, longer labels.
# ---------------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 26.7109375

# ---------------------------------------------------------------------------
# 7) Leave selection where the author left it when saving.
# ---------------------------------------------------------------------------
$ws.Range("D18").Select() | Out-Null

Write-Output "done"
